# Revert "adding term 2.0.0":
#  - Metadata sheet: Version/Date/Contact go back to the pre-2.0.0 values.
#  - "Include from FSIII" sheet: the five concept-UUID rows that were
#    inserted for 2.0.0 are removed again, restoring the original
#    J1/J5/J2/J3/J4 + System URI layout.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": revert Version, Date, Contact values ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"

# --- Sheet "Include from FSIII": drop the five UUID concept rows ---
# Rows 2-6 hold the UUID concept ids added by the reverted commit;
# deleting them shifts J1..J4 and the System URI row back up.
$wsConcepts = $wb.Worksheets.Item("Include from FSIII")
$wsConcepts.Range("A2:B6").EntireRow.Delete()
